# Remove the trailing "Ver no Jupiter / Salvar em pdf / Salvar em docx"
# paragraph and the site-footer copyright paragraph ("(c) 2020 . Contact:
# ...") from the end of the document, along with the blank paragraph that
# separated them from the "Requisitos" list above. The "LOB1038: Física
# Experimental I (Requisito fraco)" paragraph, and everything from the
# following blank paragraph onward (the page-break paragraph and section
# properties), are left untouched.

$d = $word.ActiveDocument

# Locate the "LOB1038: Física Experimental I (Requisito fraco)" paragraph
# using Find so we don't depend on a hard-coded paragraph index.
$anchor = $d.Content.Duplicate
$found = $anchor.Find.Execute(
    "LOB1038: Física Experimental I (Requisito fraco)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $lobIndex = -1
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Start -eq $anchor.Start) {
            $lobIndex = $i
            break
        }
    }

    if ($lobIndex -gt 0 -and ($lobIndex + 3) -le $d.Paragraphs.Count) {
        # The three paragraphs right after it are:
        #   lobIndex+1 -> blank paragraph
        #   lobIndex+2 -> "Ver no Jupiter Salvar em pdf Salvar em docx"
        #   lobIndex+3 -> "© 2020 . Contact: ... Creative Commons Attribution"
        $firstToDelete = $d.Paragraphs.Item($lobIndex + 1)
        $lastToDelete = $d.Paragraphs.Item($lobIndex + 3)

        $deleteRange = $d.Range($firstToDelete.Range.Start, $lastToDelete.Range.End)
        $deleteRange.Delete()
    }
}
